$wb = $excel.ActiveWorkbook

# --- AppData sheet (sheet1) ---
$appData = $wb.Worksheets.Item("AppData")

# Fill in row 5: TestName, AppURL (hyperlink), AppName, browser
$appData.Range("A5").Value = "testcase3"
$appData.Range("C5").Value = "YMO"
$appData.Range("D5").Value = "chrome"

$appData.Hyperlinks.Add($appData.Range("B5"), "http://www.yourmealsonline.co.uk")

# Column B width
$appData.Columns.Item(2).ColumnWidth = 30.140625

# --- Actions sheet (sheet2) ---
$actions = $wb.Worksheets.Item("Actions")
$actions.Range("E2").Value = "YMO"
$actions.Range("C4").Value = "no"

# Update per-sheet selections without disturbing which tab is active:
# select on Actions first, then re-select/activate AppData last so
# AppData keeps tabSelected="1" (it was the active sheet originally).
$actions.Range("H6").Select()
$appData.Range("E10").Select()
$appData.Activate()
